$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log_Muestras")

$timestamps = @{
    2  = "2025-10-17T07:09:30.715337"
    3  = "2025-10-17T07:09:30.715337"
    4  = "2025-10-17T07:09:30.715337"
    5  = "2025-10-17T07:09:30.715337"
    6  = "2025-10-17T07:09:30.715337"
    7  = "2025-10-17T07:09:30.719392"
    8  = "2025-10-17T07:09:30.719392"
    9  = "2025-10-17T07:09:30.719922"
    10 = "2025-10-17T07:09:30.719922"
    11 = "2025-10-17T07:09:30.719922"
    12 = "2025-10-17T07:09:30.719922"
    13 = "2025-10-17T07:09:30.719922"
    14 = "2025-10-17T07:09:30.719922"
    15 = "2025-10-17T07:09:30.719922"
    16 = "2025-10-17T07:09:30.790861"
    17 = "2025-10-17T07:09:30.790861"
    18 = "2025-10-17T07:09:30.791874"
    19 = "2025-10-17T07:09:30.791874"
    20 = "2025-10-17T07:09:30.791874"
    21 = "2025-10-17T07:09:30.792842"
    22 = "2025-10-17T07:09:30.793261"
    23 = "2025-10-17T07:09:30.793261"
    24 = "2025-10-17T07:09:30.793261"
    25 = "2025-10-17T07:09:30.793261"
    26 = "2025-10-17T07:09:30.860714"
    27 = "2025-10-17T07:09:30.860714"
    28 = "2025-10-17T07:09:30.860714"
    29 = "2025-10-17T07:09:30.861713"
    30 = "2025-10-17T07:09:30.861713"
    31 = "2025-10-17T07:09:30.861713"
    32 = "2025-10-17T07:09:30.861713"
    33 = "2025-10-17T07:09:30.861713"
    34 = "2025-10-17T07:09:30.862714"
    35 = "2025-10-17T07:09:30.862714"
    36 = "2025-10-17T07:09:30.862714"
    37 = "2025-10-17T07:09:30.862714"
    38 = "2025-10-17T07:09:30.862714"
    39 = "2025-10-17T07:09:30.863718"
    40 = "2025-10-17T07:09:30.863718"
    41 = "2025-10-17T07:09:30.863718"
    42 = "2025-10-17T07:09:30.863718"
    43 = "2025-10-17T07:09:30.863718"
    44 = "2025-10-17T07:09:30.864715"
    45 = "2025-10-17T07:09:30.864715"
    46 = "2025-10-17T07:09:30.864715"
    47 = "2025-10-17T07:09:30.864715"
    48 = "2025-10-17T07:09:30.864715"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 26).Value = $timestamps[$row]
}
